$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function SetTextValue($Cell, $Value) {
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.Style = "Normal"
}

# Row 2
SetTextValue $ws.Range("D2") '42.972.79'
SetTextValue $ws.Range("E2") '  +2.11%  '

# Row 3
SetTextValue $ws.Range("D3") '2.298.63'
SetTextValue $ws.Range("E3") '  +1.57%  '

# Row 4
SetTextValue $ws.Range("E4") '  -0.01%  '

# Row 5
SetTextValue $ws.Range("D5") '301.80'
SetTextValue $ws.Range("E5") '  +0.76%  '

# Row 6
SetTextValue $ws.Range("D6") '98.94'
SetTextValue $ws.Range("E6") '  +5.00%  '

# Row 7
SetTextValue $ws.Range("D7") '0.506'
SetTextValue $ws.Range("E7") '  +2.04%  '

# Row 8
SetTextValue $ws.Range("E8") '  -0.05%  '

# Row 9
SetTextValue $ws.Range("E9") '  +3.10%  '

# Row 10
SetTextValue $ws.Range("D10") '34.19'
SetTextValue $ws.Range("E10") '  +3.53%  '

# Row 11
SetTextValue $ws.Range("E11") '  +1.12%  '

# Row 12
SetTextValue $ws.Range("D12") '49.14'
SetTextValue $ws.Range("E12") '  +2.59%  '

# Row 13
SetTextValue $ws.Range("E13") '  +4.20%  '

# Row 14
SetTextValue $ws.Range("D14") '17.81'
SetTextValue $ws.Range("E14") '  +15.63%  '

# Row 15
SetTextValue $ws.Range("D15") '6.78'
SetTextValue $ws.Range("E15") '  +1.53%  '

# Row 16
SetTextValue $ws.Range("D16") '2.655.35'
SetTextValue $ws.Range("E16") '  +1.56%  '

# Row 17
SetTextValue $ws.Range("D17") '2.258.90'
SetTextValue $ws.Range("E17") '  -0.19%  '

# Row 18
SetTextValue $ws.Range("D18") '0.807'
SetTextValue $ws.Range("E18") '  +4.46%  '

# Row 19
SetTextValue $ws.Range("D19") '42.878.37'
SetTextValue $ws.Range("E19") '  +1.91%  '

# Row 20
SetTextValue $ws.Range("D20") '12.36'
SetTextValue $ws.Range("E20") '  +8.87%  '

# Row 21
SetTextValue $ws.Range("D21") '0.0₃0902'
SetTextValue $ws.Range("E21") '  +1.23%  '

# Row 22
SetTextValue $ws.Range("E22") '  +1.63%  '

# Row 23
SetTextValue $ws.Range("D23") '67.77'
SetTextValue $ws.Range("E23") '  +1.68%  '

# Row 24
SetTextValue $ws.Range("D24") '236.17'
SetTextValue $ws.Range("E24") '  +1.19%  '

# Row 25
SetTextValue $ws.Range("E25") '  +13.77%  '

# Row 26
SetTextValue $ws.Range("D26") '0.999'
SetTextValue $ws.Range("E26") '  -0.13%  '

# Row 27
SetTextValue $ws.Range("E27") '  -0.01%  '

# Row 28
SetTextValue $ws.Range("D28") '24.71'
SetTextValue $ws.Range("E28") '  +3.99%  '

# Row 29
SetTextValue $ws.Range("B29") 'Toncoin'
SetTextValue $ws.Range("C29") 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
SetTextValue $ws.Range("D29") '2.18'
SetTextValue $ws.Range("E29") '  -3.50%  '

# Row 30
SetTextValue $ws.Range("B30") 'Monero'
SetTextValue $ws.Range("C30") 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
SetTextValue $ws.Range("D30") '167.96'
SetTextValue $ws.Range("E30") '  +0.48%  '

# Row 31
SetTextValue $ws.Range("B31") 'Cosmos'
SetTextValue $ws.Range("C31") 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
SetTextValue $ws.Range("D31") '9.15'
SetTextValue $ws.Range("E31") '  +1.16%  '

# Row 32
SetTextValue $ws.Range("B32") 'InjectiveProtocol'
SetTextValue $ws.Range("C32") 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
SetTextValue $ws.Range("D32") '33.57'
SetTextValue $ws.Range("E32") '  -0.58%  '

# Row 33
SetTextValue $ws.Range("E33") '  +0.02%  '

# Row 34
SetTextValue $ws.Range("D34") '5.04'
SetTextValue $ws.Range("E34") '  +2.12%  '

# Row 35
SetTextValue $ws.Range("B35") 'WEMIXToken'
SetTextValue $ws.Range("C35") 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
SetTextValue $ws.Range("D35") '2.43'
SetTextValue $ws.Range("E35") '  +3.68%  '

# Row 36
SetTextValue $ws.Range("B36") 'RenderToken'
SetTextValue $ws.Range("C36") 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
SetTextValue $ws.Range("D36") '4.54'
SetTextValue $ws.Range("E36") '  +0.98%  '

# Row 37
SetTextValue $ws.Range("D37") '16.86'
SetTextValue $ws.Range("E37") '  +4.44%  '

# Row 38
SetTextValue $ws.Range("D38") '0.0692'
SetTextValue $ws.Range("E38") '  -0.16%  '

# Row 39
SetTextValue $ws.Range("E39") '  +3.04%  '

# Row 40
SetTextValue $ws.Range("E40") '  +4.51%  '

# Row 41
SetTextValue $ws.Range("E41") '  +0.38%  '

# Row 42
SetTextValue $ws.Range("E42") '  -0.01%  '

# Row 43
SetTextValue $ws.Range("D43") '2.37'
SetTextValue $ws.Range("E43") '  -1.89%  '

# Row 44
SetTextValue $ws.Range("D44") '1.994.13'
SetTextValue $ws.Range("E44") '  +1.95%  '

# Row 45
SetTextValue $ws.Range("D45") '0.0285'
SetTextValue $ws.Range("E45") '  +2.15%  '

# Row 46
SetTextValue $ws.Range("E46") '  +4.93%  '

# Row 47
SetTextValue $ws.Range("D47") '17.53'
SetTextValue $ws.Range("E47") '  +0.54%  '

# Row 48
SetTextValue $ws.Range("E48") '  +2.18%  '

# Row 49
SetTextValue $ws.Range("D49") '56.88'
SetTextValue $ws.Range("E49") '  +9.91%  '

# Row 50
SetTextValue $ws.Range("D50") '2.526.76'
SetTextValue $ws.Range("E50") '  +1.55%  '

# Row 51
SetTextValue $ws.Range("E51") '  +3.30%  '
